$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E10 value (855528 -> 854877)
$ws.Range("E10").Value = 854877

# Add new row 11 data, copying the formatting from A10 (style "s=1")
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = 9

$ws.Range("B11").Value = "M2_10 Cat 2020"

$ws.Range("C11").Value = 9703
$ws.Range("D11").Value = 10804
$ws.Range("E11").Value = 929613
$ws.Range("F11").Value = 9977
$ws.Range("G11").Value = 10067
$ws.Range("H11").Value = 10176
